# bug fixed / v.4.0
# Append new incident rows (123-126) to the activity log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("WC48 P5F", "Cámara no detecta skeleton", "'2024-06-04", "'09:12:54", "Mañana", "'09:12:56", "'0:00:02", "N/A"),
    @("WC48 P5F", "Etiquetadora",               "'2024-06-04", "'09:13:23", "Mañana", "'09:13:25", "'0:00:02", "N/A"),
    @("WC48 P5F", "AOI (fallo etiqueta)",        "'2024-06-04", "'09:13:28", "Mañana", "'09:13:29", "'0:00:01", "0.49 minutos"),
    @("WC49 P5H", "Tornillo atascado",           "'2024-06-04", "'09:16:44", "Mañana", "'09:16:44", "'0:00:00", "N/A")
)

$startRow = 123
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
